$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 5 new variable rows (error, decision_error, one_tailed,
#        error_1tail, decision_error_1tail) right before the existing
#        "extract_apa" row (current row 31). ---
$ws.Rows("31:35").Insert()

# Copy the formatting (borders etc.) of the row directly above down into
# the freshly inserted rows so they look like the rest of the table.
$ws.Range("A30:C30").Copy()
$ws.Range("A31:C35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A31").Value = "error"
$ws.Range("B31").Value = "p-value does not match degrees of freedom and test statistic"
$ws.Range("C31").Value = "logical"

$ws.Range("A32").Value = "decision_error"
$ws.Range("B32").Value = "recomputed p-value is significant and reported is not, or vice versa"
$ws.Range("C32").Value = "logical"

$ws.Range("A33").Value = "one_tailed"
$ws.Range("B33").Value = "is this test explicitly identified as a one-tailed test in the text?"
$ws.Range("C33").Value = "logical"

$ws.Range("A34").Value = "error_1tail"
$ws.Range("B34").Value = "is the result an error when taking into account one-tailed testing?"
$ws.Range("C34").Value = "logical"

$ws.Range("A35").Value = "decision_error_1tail"
$ws.Range("B35").Value = "is the result a decision error when taking into account one-tailed testing?"
$ws.Range("C35").Value = "logical"

# --- 2. The previously existing "extract_apa" and "pdf_conversion_issues"
#        rows are now rows 36 and 37; their "type" column changes from
#        "factor" to "logical". ---
$ws.Range("C36").Value = "logical"
$ws.Range("C37").Value = "logical"

# --- 3. Insert a new "typesetting_issues" row right after
#        "pdf_conversion_issues" (before "remarks", currently row 38). ---
$ws.Rows("38:38").Insert()
$ws.Range("A37:C37").Copy()
$ws.Range("A38:C38").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A38").Value = "typesetting_issues"
$ws.Range("B38").Value = "are there typesetting issues? (e.g., result spanned two pages)"
$ws.Range("C38").Value = "logical"

# --- 4. The explanation (column B) cells across the whole variable table
#        lose their fill/border styling (only the outer box border on
#        columns A and C remains) - clear the formatting on column B. ---
$ws.Range("B13:B38").ClearFormats()
$ws.Range("B13").ClearContents()
$ws.Range("B14").ClearContents()

# --- 5. Update the view: the editor ended up with C37 selected and the
#        window scrolled down so row 9 is the first visible row. ---
$ws.Range("C37").Select()
$excel.ActiveWindow.ScrollRow = 9

Write-Output "edit complete"
